# Update the Entities sheet so that Item_3 / Item_4 become the new puzzle
# block items, sharing one localized description, and move the active
# selection to C5 (matching the saved state captured in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Entities")

# Rename the two puzzle items first so their shared-string entries are
# interned before the (shared) description text.
$ws.Range("B4").Value = "Item_PuzzleBlock_A"
$ws.Range("B5").Value = "Item_PuzzleBlock_B"

# Row 4 (id=3) and row 5 (id=4) both use the same localized description
$ws.Range("C4").Value = "パズルブロッカだ。パズル用みたいです。"
$ws.Range("C5").Value = "パズルブロッカだ。パズル用みたいです。"

# Move/save the active cell selection as it was left in the authored file
$ws.Activate()
$ws.Range("C5").Select()

$wb.Save()
